$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2: new "Button" header in E1 (plain, no special formatting)
$ws2.Range("E1").Value = "Button"

# Sheet1: new "Submit" header in E1, styled like the other header cells
# (bold, centered/top aligned, thin left/right border, no top/bottom border)
$ws1.Range("E1").Value = "Submit"
$ws1.Range("E1").Font.Bold = $true
$ws1.Range("E1").HorizontalAlignment = -4108
$ws1.Range("E1").VerticalAlignment = -4160
$ws1.Range("E1").Interior.ColorIndex = -4142
$ws1.Range("E1").Borders.LineStyle = 1
$ws1.Range("E1").Borders.Item(8).LineStyle = -4142
$ws1.Range("E1").Borders.Item(9).LineStyle = -4142

# Move the selection in each sheet to the new header cell, matching the
# post-edit Excel state (Sheet1 stays the active tab).
$ws2.Range("E1").Select()
$ws1.Range("E1").Select()
